# Update odds values on Sheet1 as per the latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (San Diego FC vs Minnesota Utd)
$ws.Range("K2").Value = 4.8
$ws.Range("U2").Value = 2.2

# Row 4 (Vitesse Arnhem vs Jong PSV Eindhoven)
$ws.Range("Q4").Value = 1.42

# Row 6 (Roda JC vs FC Dordrecht)
$ws.Range("G6").Value = 2.14
$ws.Range("H6").Value = 3.45
$ws.Range("I6").Value = 4.4
$ws.Range("J6").Value = 3.45

# Row 7 (RKC Waalwijk vs MVV Maastricht)
$ws.Range("F7").Value = 1.45
$ws.Range("G7").Value = 1.53
$ws.Range("J7").Value = 5.1
